$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("models")

# Remove the trailing "[1]" citation marker from the ARIMA description in C2.
$old = $ws.Range("C2").Value2
$new = $old -replace '\[1\]$', ''
$ws.Range("C2").Value = $new

# Update the saved selection to match the authored state (C8).
[void]$ws.Range("C8").Select()
